# Fruta / hortaliza, semanal
# A new weekly price record (2021-11-26) is inserted as a new data row
# right after the header/first rows, pushing the existing history rows
# down by one (old row 5 -> row 6, old row 6 -> row 7, ..., old row 15 -> row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5 - shifts rows 5..15 down to 6..16
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the latest week's record
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C5").Value = 'Los Lagos'
$ws.Range("D5").Value = 44526
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = 'Otros'
$ws.Range("I5").Value = 100107002
$ws.Range("J5").Value = 'Chirimoya'
$ws.Range("K5").Value = 'Cultivar IV Región'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 21000
$ws.Range("Q5").Value = '$/bandeja 8 kilos'
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 2625
$ws.Range("T5").Value = 8
